# Update the Q&A table: questions/answers were reworded and reordered,
# a couple of brand-new rows' worth of content rotated in, and the
# column widths / row heights were touched up to fit the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell contents -------------------------------------------------------
$ws.Range("A1").Value = "Pregunta"
$ws.Range("B1").Value = "Respuesta"

$ws.Range("A2").Value = "¿Se puede cerrar una sucursal en horario de atención a público?"
$ws.Range("B2").Value = "No, no se puede cerrar una surcursal, a no ser que existan razones extraordinarias que justifiquen su cierre."

$ws.Range("A3").Value = "¿Cuándo se puede cerrar una surcusal en horario de atención de público?"
$ws.Range("B3").Value = "Se puede cerrar una sucursal si existen razones extraordinarias como emergencias sanitarias, cortes de electricidad y agua prolongados, etc. Que justifique el cierre"

$ws.Range("A4").Value = "¿Qué se debe hacer si por motivos extraordinarios se `ndebe cerrar una sucursal?"
$ws.Range("B4").Value = "Ante un cierre por emergencia u programado de una sucursal se deben tomar los resguardo necesarios para darle continuidad al servicio y asegurar la atención a las y los trabajadores.`nEl DR o su subrogante debe contactar al Departamento de Atención de Usuarios , comunicar los motivos del cierre y planificar la estrategia que asegure la atención oportuna al público que lo requiera."

$ws.Range("A5").Value = "La o el Seremi de mi Región solicitó un vehículo institucional para una actividad de la Seremía.`n ¿Se puede `"prestar`" el vehículo?"
$ws.Range("B5").Value = "No, los vehículos que se utilizan en el ISL no son vehículos fiscales, son arrendados y no pueden ser utilizados en ninguna función que no`nesté contemplada en el contrato de arriendo respectivo."

$ws.Range("A6").Value = "La Directora o Director Nacional viene a la Región, llega al aeropuerto y e hará uso del vehículo institucional`npara facilitar su traslado."
$ws.Range("B6").Value = "No, los vehículos que se utilizan en el ISL no son vehículos fiscales, son arrendados y no pueden ser utilizados en ninguna función que no`nesté contemplada en el contrato de arriendo respectivo."

$ws.Range("A7").Value = "Existe un problema con los estacionamientos y los vehículos no pueden quedar guardados después del horario laboral.`n¿Puede el funcionario a cargo llevarse el vehículo a su hogar, para que quede resguardado?"
$ws.Range("B7").Value = "No, ningún funcionario o funcionaria puede llevarse el vehículo institucional a su hogar, independiente del motivo."

$ws.Range("A8").Value = "Se cortó la luz. `n¿Es posible activar las VPN para que las o los funcionarios trabajen desde sus hogares?"

# --- Row heights: rows 6 & 7 now hold shorter text than before ----------
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 60

# --- Column widths widened slightly to fit the revised wording ----------
$ws.Columns.Item(1).ColumnWidth = 67.666667
$ws.Columns.Item(2).ColumnWidth = 144.5

# --- Selection moved to A9 after the edits -------------------------------
[void]$ws.Range("A9").Select()
